$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D. Excel shifts old D:K data to F:M.
$ws.Columns("D:E").Insert(-4161, 0)

# The two new columns copied formatting (s=1) from column C on insert; restore the
# correct number/date formatting by pasting it across from column F, in the three
# contiguous blocks of data rows (skipping the blank/separator rows 36-37 & 78-79).
$ws.Range("F7:F35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Range("F7:F35").Copy()
$ws.Range("E7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("E38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("E80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarters of data, and apply the handful of revised historical values
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 2230800
$ws.Range("E8").Value = 2192500
$ws.Range("D9").Value = 1831900
$ws.Range("E9").Value = 1805800
$ws.Range("D10").Value = 398900
$ws.Range("E10").Value = 386700
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 13500
$ws.Range("H14").Value = -12400
$ws.Range("D15").Value = 12000
$ws.Range("E15").Value = 13100
$ws.Range("D17").Value = 2205600
$ws.Range("E17").Value = 2177000
$ws.Range("H17").Value = 2190900
$ws.Range("D18").Value = 25200
$ws.Range("E18").Value = 15500
$ws.Range("H18").Value = 33500
$ws.Range("D20").Value = 2500
$ws.Range("E20").Value = 400
$ws.Range("H20").Value = 12400
$ws.Range("D21").Value = 39700
$ws.Range("E21").Value = 29000
$ws.Range("D22").Value = 11800
$ws.Range("E22").Value = 11000
$ws.Range("D23").Value = 15900
$ws.Range("E23").Value = 4900
$ws.Range("D24").Value = -24900
$ws.Range("E24").Value = 3500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 40800
$ws.Range("E26").Value = 1400
$ws.Range("D27").Value = 40800
$ws.Range("E27").Value = 1400
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -31500
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -2500
$ws.Range("E32").Value = -400
$ws.Range("H32").Value = -12400
$ws.Range("D33").Value = 9300
$ws.Range("E33").Value = 1400
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 9300
$ws.Range("E35").Value = 1400
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 64300
$ws.Range("E41").Value = 70900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 1184600
$ws.Range("E43").Value = 1215700
$ws.Range("D44").Value = 688200
$ws.Range("E44").Value = 737800
$ws.Range("D45").Value = 147200
$ws.Range("E45").Value = 152500
$ws.Range("D46").Value = 2084300
$ws.Range("E46").Value = 2176900
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 206700
$ws.Range("E48").Value = 210100
$ws.Range("D49").Value = 156800
$ws.Range("E49").Value = 158100
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 81900
$ws.Range("E52").Value = 89400
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 2529700
$ws.Range("E54").Value = 2634500
$ws.Range("D57").Value = 641900
$ws.Range("E57").Value = 723400
$ws.Range("D58").Value = 7300
$ws.Range("E58").Value = 7300
$ws.Range("D59").Value = 200500
$ws.Range("E59").Value = 207300
$ws.Range("D60").Value = 849700
$ws.Range("E60").Value = 938000
$ws.Range("D61").Value = 987200
$ws.Range("E61").Value = 1022700
$ws.Range("D62").Value = 149700
$ws.Range("E62").Value = 137100
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1986600
$ws.Range("E66").Value = 2097800
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -8500
$ws.Range("E72").Value = -17800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 543100
$ws.Range("E76").Value = 536700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 9300
$ws.Range("E81").Value = 1400
$ws.Range("D83").Value = 12000
$ws.Range("E83").Value = 13100
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 7500
$ws.Range("E89").Value = -500
$ws.Range("D91").Value = -11700
$ws.Range("E91").Value = -12200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 7900
$ws.Range("E94").Value = -12100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -21700
$ws.Range("E100").Value = 13900
$ws.Range("D101").Value = -300
$ws.Range("E101").Value = 100
$ws.Range("D102").Value = -6600
$ws.Range("E102").Value = 1400
